$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hendrik")
$ws.Activate()

# New rows of data: date serial, start time, end time, total hours, description
$rows = @(
    @{ Date = 44726; Start = 0.625;              End = 0.70833333333333337; Total = 2;   Desc = "Overlopen en verstaan code" },
    @{ Date = 44728; Start = 0.75;                End = 0.77083333333333337; Total = 0.5; Desc = "prestantie" },
    @{ Date = 44729; Start = 0.58333333333333337; End = 0.66666666666666663; Total = 2;   Desc = "Class Diagram" },
    @{ Date = 44730; Start = 0.70833333333333337; End = 0.79166666666666663; Total = 2;   Desc = "Opmaak bijkomstige info" }
)

$r = 17
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Date
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat

    $ws.Cells.Item($r, 2).Value = $row.Start
    $ws.Cells.Item($r, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat

    $ws.Cells.Item($r, 3).Value = $row.End
    $ws.Cells.Item($r, 3).NumberFormat = $ws.Cells.Item(16, 3).NumberFormat

    $ws.Cells.Item($r, 4).Value = $row.Total

    $ws.Cells.Item($r, 5).Value = $row.Desc

    $r++
}

$ws.Range("E20").Select()
